# Updated cryptos list on Sun Jan 28 14:29:17 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    # Force the cell to remain plain text even if the string looks
    # like a number (e.g. "1.00", "36.10"), so trailing zeros /
    # formatting are preserved exactly as in the source data.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" "42.464.72"
$ws.Range("E2").Value = "  +1.64%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.284.72"
$ws.Range("E3").Value = "  +0.65%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "307.74"
$ws.Range("E5").Value = "  +1.47%  "

# Row 6 - Solana
Set-TextValue "D6" "98.11"
$ws.Range("E6").Value = "  +5.87%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.01%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.496"
$ws.Range("E9").Value = "  +2.44%  "

# Row 10 - Avalanche
Set-TextValue "D10" "36.10"
$ws.Range("E10").Value = "  +10.75%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +0.31%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -0.93%  "

# Row 13 - Polkadot
Set-TextValue "D13" "6.73"
$ws.Range("E13").Value = "  +0.84%  "

# Row 14 - Wrapped liquid staked Ether 2.0
Set-TextValue "D14" "2.637.63"
$ws.Range("E14").Value = "  +0.61%  "

# Row 15 - Chainlink
Set-TextValue "D15" "14.49"
$ws.Range("E15").Value = "  +1.51%  "

# Row 16 - Wrapped Ether
Set-TextValue "D16" "2.297.70"
$ws.Range("E16").Value = "  +1.05%  "

# Row 17 - Polygon
Set-TextValue "D17" "0.799"
$ws.Range("E17").Value = "  +3.10%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "42.369.99"
$ws.Range("E18").Value = "  +1.66%  "

# Row 19 - Internet Computer (DFINITY)
Set-TextValue "D19" "12.60"
$ws.Range("E19").Value = "  +0.88%  "

# Row 20 - Shiba Inu
Set-TextValue "D20" "0.0₃0914"
$ws.Range("E20").Value = "  +1.02%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +0.80%  "

# Row 22 - Litecoin
Set-TextValue "D22" "67.87"
$ws.Range("E22").Value = "  +1.20%  "

# Row 23 - Bitcoin Cash
Set-TextValue "D23" "241.38"
$ws.Range("E23").Value = "  +0.58%  "

# Row 24 - PancakeSwap
Set-TextValue "D24" "2.61"
$ws.Range("E24").Value = "  +1.33%  "

# Row 25 - ImmutableX
$ws.Range("E25").Value = "  +1.67%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.02%  "

# Row 27 - Ethereum Classic
Set-TextValue "D27" "23.96"
$ws.Range("E27").Value = "  -0.12%  "

# Row 28 - Injective Protocol
Set-TextValue "D28" "37.87"

# Row 29 - Cosmos
$ws.Range("E29").Value = "  +0.20%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  +1.64%  "

# Row 31 - Monero
Set-TextValue "D31" "159.50"
$ws.Range("E31").Value = "  -0.59%  "

# Row 32 - Filecoin
Set-TextValue "D32" "5.27"
$ws.Range("E32").Value = "  +0.45%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  +4.75%  "

# Row 35 - Hedera
Set-TextValue "D35" "0.0745"
$ws.Range("E35").Value = "  +0.09%  "

# Row 36 - Celestia
$ws.Range("E36").Value = "  +1.27%  "

# Row 37 - Kaspa
$ws.Range("E37").Value = "  +0.92%  "

# Row 38 - WEMIX Token
$ws.Range("E38").Value = "  +0.09%  "

# Row 39 - ARBITRUM
Set-TextValue "D39" "1.85"
$ws.Range("E39").Value = "  +3.18%  "

# Row 40 - Stellar
$ws.Range("E40").Value = "  -1.10%  "

# Row 41 - RenderToken
$ws.Range("E41").Value = "  +5.53%  "

# Row 42 - ApeXProtocol
$ws.Range("E42").Value = "  +15.01%  "

# Row 43 - Maker
Set-TextValue "D43" "2.002.11"
$ws.Range("E43").Value = "  -0.05%  "

# Row 44 - VeChain
Set-TextValue "D44" "0.0286"
$ws.Range("E44").Value = "  +1.89%  "

# Row 45 - EnergySwap
Set-TextValue "D45" "18.95"
$ws.Range("E45").Value = "  -1.05%  "

# Row 46 - NEAR Protocol
Set-TextValue "D46" "2.99"
$ws.Range("E46").Value = "  +3.37%  "

# Row 47 - FraxShare
$ws.Range("E47").Value = "  -3.09%  "

# Row 48 - MultiversX
Set-TextValue "D48" "53.12"
$ws.Range("E48").Value = "  +1.09%  "

# Row 49 - Stacks
$ws.Range("E49").Value = "  +1.35%  "

# Row 50 - BitcoinSV
Set-TextValue "D50" "72.27"

# Row 51 - Aave
Set-TextValue "D51" "92.43"
